$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-8 from serial 45233 to 45243
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45243
}
